# Generate Report for Handoff
#
# Adds two new files that are "Ready for handoff" to the localization status
# report:
#   5344c2e5-07af-494c-99ae-8305e3c2db96.md
#   c59f5440-6967-45b7-9310-88f54c37ed63.md
#
# They are inserted right before the always-last ".localization-config" row
# on every worksheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/437e2241346673991de55533f87aa42d6066ef5b/e2e/"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/437e2241346673991de55533f87aa42d6066ef5b/.localization-config"
$zhcnBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/462cc301c34512973d9d81f633aef88009c7d1d5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/"
$dedeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/932f812ae22d67118469813f289dee875c572802/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/"

$uuid1 = "5344c2e5-07af-494c-99ae-8305e3c2db96"
$hash1 = "4a43194d10cf3d5a29e23050d5c64ba8f9ae4670"
$uuid2 = "c59f5440-6967-45b7-9310-88f54c37ed63"
$hash2 = "054b334e3af95f8bb45407edc9a6abac2852311e"

$statusText = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn, C=de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Push the ".localization-config" row (row 4) down by two, making room for
# the two new "ready for handoff" rows.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4,2).Value = $statusText
$ws.Cells.Item(4,3).Value = $statusText

$ws.Cells.Item(5,2).Value = $statusText
$ws.Cells.Item(5,3).Value = $statusText

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdBase + "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.md", "", "", "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $mdBase + "f8d3bd2b-b640-4983-a712-a08aaa377cb4.md", "", "", "f8d3bd2b-b640-4983-a712-a08aaa377cb4.md")
$ws.Hyperlinks.Add($ws.Range("A4"), $mdBase + $uuid1 + ".md", "", "", $uuid1 + ".md")
$ws.Hyperlinks.Add($ws.Range("A5"), $mdBase + $uuid2 + ".md", "", "", $uuid2 + ".md")
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": detail table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4,2).Value = $statusText
$ws.Cells.Item(4,4).Value = "2016-01-25 07:31:38"
$ws.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4,8).Value = "Include"

$ws.Cells.Item(5,2).Value = $statusText
$ws.Cells.Item(5,4).Value = "2016-01-25 07:31:38"
$ws.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(5,8).Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdBase + "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.md", "", "", "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $zhcnBase + "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.7466a6110f46d0932efe872f1d4d26637330b416.zh-cn.xlf", "", "", "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.7466a6110f46d0932efe872f1d4d26637330b416.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $mdBase + "f8d3bd2b-b640-4983-a712-a08aaa377cb4.md", "", "", "f8d3bd2b-b640-4983-a712-a08aaa377cb4.md")
$ws.Hyperlinks.Add($ws.Range("C3"), $zhcnBase + "f8d3bd2b-b640-4983-a712-a08aaa377cb4.a1384396d74f7b81c06ebbed417a703f69d932c9.zh-cn.xlf", "", "", "f8d3bd2b-b640-4983-a712-a08aaa377cb4.a1384396d74f7b81c06ebbed417a703f69d932c9.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), $mdBase + $uuid1 + ".md", "", "", $uuid1 + ".md")
$ws.Hyperlinks.Add($ws.Range("C4"), $zhcnBase + $uuid1 + "." + $hash1 + ".zh-cn.xlf", "", "", $uuid1 + "." + $hash1 + ".zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), $mdBase + $uuid2 + ".md", "", "", $uuid2 + ".md")
$ws.Hyperlinks.Add($ws.Range("C5"), $zhcnBase + $uuid2 + "." + $hash2 + ".zh-cn.xlf", "", "", $uuid2 + "." + $hash2 + ".zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de": detail table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4,2).Value = $statusText
$ws.Cells.Item(4,4).Value = "2016-01-25 07:31:50"
$ws.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4,8).Value = "Include"

$ws.Cells.Item(5,2).Value = $statusText
$ws.Cells.Item(5,4).Value = "2016-01-25 07:31:50"
$ws.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(5,8).Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdBase + "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.md", "", "", "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $dedeBase + "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.7466a6110f46d0932efe872f1d4d26637330b416.de-de.xlf", "", "", "7da0ee49-e0a5-4455-8fb6-a8b921b9df64.7466a6110f46d0932efe872f1d4d26637330b416.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $mdBase + "f8d3bd2b-b640-4983-a712-a08aaa377cb4.md", "", "", "f8d3bd2b-b640-4983-a712-a08aaa377cb4.md")
$ws.Hyperlinks.Add($ws.Range("C3"), $dedeBase + "f8d3bd2b-b640-4983-a712-a08aaa377cb4.a1384396d74f7b81c06ebbed417a703f69d932c9.de-de.xlf", "", "", "f8d3bd2b-b640-4983-a712-a08aaa377cb4.a1384396d74f7b81c06ebbed417a703f69d932c9.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), $mdBase + $uuid1 + ".md", "", "", $uuid1 + ".md")
$ws.Hyperlinks.Add($ws.Range("C4"), $dedeBase + $uuid1 + "." + $hash1 + ".de-de.xlf", "", "", $uuid1 + "." + $hash1 + ".de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), $mdBase + $uuid2 + ".md", "", "", $uuid2 + ".md")
$ws.Hyperlinks.Add($ws.Range("C5"), $dedeBase + $uuid2 + "." + $hash2 + ".de-de.xlf", "", "", $uuid2 + "." + $hash2 + ".de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config")

$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()
